$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gyn")

# --- Header row: duplicate the "Gyn" header into new column D ---
$ws.Range("D1").Value = "Gyn"

# --- Row 2: Encounter Date - update existing date, add new column D date ---
$ws.Range("C2").Value = 43899
$ws.Range("D2").Value = 43899
$ws.Range("D2").NumberFormat = "yyyy\-mm\-dd;@"

# --- Row 3: Responsibility / Observed ---
$ws.Range("D3").Value = 1

# --- Row 6: Age ---
$ws.Range("C6").Value = 47
$ws.Range("D6").Value = 54

# --- Row 9: Setting / OR ---
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1

# --- Row 10: Setting / Outpatient - clear old mark ---
$ws.Range("C10").ClearContents()

# --- Row 12: Diagnoses / Abnormal uterine bleeding ---
$ws.Range("C12").ClearContents()
$ws.Range("D12").Value = 1

# --- Row 13: Acute Vaginitis and STIs - clear old mark ---
$ws.Range("C13").ClearContents()

# --- Row 15: Benign pelvic and adnexal masses ---
$ws.Range("D15").Value = 1

# --- Row 22: Management of abnormal pap smear ---
$ws.Range("D22").Value = 1

# --- Row 33: Bimanual pelvic palpation ---
$ws.Range("D33").Value = 1

# --- Row 36: Obtaining pap smear - clear old mark ---
$ws.Range("C36").ClearContents()

# --- Row 37: Speculum examination - clear old mark ---
$ws.Range("C37").ClearContents()

# --- Row 41: Vaginal and cervical culture acquisition - clear old mark ---
$ws.Range("C41").ClearContents()

# --- Row 44: Performed history - clear old mark ---
$ws.Range("C44").ClearContents()

# --- Row 45: Performed Clinical Exam ---
$ws.Range("D45").Value = 1

# --- Row 47: Wrote patient note - clear old mark ---
$ws.Range("C47").ClearContents()

# --- Row 48: Performed oral patient presentation - clear old mark ---
$ws.Range("C48").ClearContents()

# --- Row 54: Gynecologic pelvic ultrasound - clear old mark ---
$ws.Range("C54").ClearContents()

# --- Row 62: new "DONE?" marker row ---
$ws.Range("A62").Value = "DONE?"
$ws.Range("B62").Value = "DONE?"
$ws.Range("B62").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("C62").Value = 1
$ws.Range("D62").Value = 1

# --- Column widths: make column D match column C ---
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# --- View state: scroll down and select E62 like the saved workbook ---
$ws.Activate()
$ws.Range("E62").Select()
